# feat: calculate validation step
# Update the AVG (D) and STD (E) columns of the stats table on Sheet1
# with the recalculated values produced by the new validation step.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = 44.76674937965261
    "E2"  = 17.66338614348757
    "D3"  = 164.5359801488834
    "E3"  = 12.51985768836525
    "D4"  = 80.09454094292805
    "E4"  = 22.73292087518887
    "D5"  = 29.42257563078581
    "E5"  = 10.17351018593594
    "D6"  = 95.71712158808933
    "E6"  = 18.25603913307767
    "D7"  = 36.52357320099256
    "E7"  = 4.441584937014361
    "D8"  = 0.6923076923076923
    "E8"  = 0.810144256311929
    "D9"  = 1.75682382133995
    "E9"  = 1.798402230679146
    "D10" = 6.173697270471464
    "E10" = 0.6762517250913661
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
